$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 296
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 7353.5713
$ws.Range("I33").Value = 128.57143
$ws.Range("J33").Value = 14578.571
$ws.Range("K33").Value = 128.57143
$ws.Range("L33").Value = 14578.571
$ws.Range("M33").Value = 100.42857
$ws.Range("N33").Value = -15036.571

# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 1186.591
$ws.Range("I70").Value = 975
$ws.Range("J70").Value = 1398.1818
$ws.Range("K70").Value = 2925
$ws.Range("L70").Value = 4194.5454
$ws.Range("M70").Value = -2655
$ws.Range("N70").Value = -4734.5454

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 1186.591
$ws.Range("I73").Value = 975
$ws.Range("J73").Value = 1398.1818
$ws.Range("K73").Value = 2925
$ws.Range("L73").Value = 4194.5454
$ws.Range("M73").Value = -1989
$ws.Range("N73").Value = -6066.5454

# Row 98: The Dotted Line
$ws.Range("H98").Value = 11121257
$ws.Range("I98").Value = 12266
$ws.Range("J98").Value = 55557224
$ws.Range("K98").Value = 12266
$ws.Range("L98").Value = 55557224
$ws.Range("M98").Value = -10768
$ws.Range("N98").Value = -55560220

# Row 122: Wishful Inking
$ws.Range("H122").Value = 11121257
$ws.Range("I122").Value = 12266
$ws.Range("J122").Value = 55557224
$ws.Range("K122").Value = 36798
$ws.Range("L122").Value = 166671672
$ws.Range("M122").Value = -34348
$ws.Range("N122").Value = -166676572

# Row 125: Body over Mind
$ws.Range("H125").Value = 2175
$ws.Range("J125").Value = 2550
$ws.Range("L125").Value = 22950
$ws.Range("N125").Value = -27870

# Row 129: Practical Command
$ws.Range("H129").Value = 1015.4815
$ws.Range("I129").Value = 255.125
$ws.Range("J129").Value = 1147.7174
$ws.Range("K129").Value = 765.375
$ws.Range("L129").Value = 3443.1522
$ws.Range("M129").Value = 4234.625
$ws.Range("N129").Value = -13443.1522

# Row 135: For Tired Minds
$ws.Range("H135").Value = 5421.8
$ws.Range("I135").Value = 3529.1
$ws.Range("J135").Value = 9207.200000000001
$ws.Range("K135").Value = 31761.9
$ws.Range("L135").Value = 82864.8
$ws.Range("M135").Value = -29226.9
$ws.Range("N135").Value = -87934.8

# Row 136: I Like Big Brush and I Cannot Lie
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 40000
$ws.Range("N136").Value = -50200

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1167.6471
$ws.Range("I137").Value = 841.7727
$ws.Range("J137").Value = 1765.0834
$ws.Range("K137").Value = 2525.3181
$ws.Range("L137").Value = 5295.2502
$ws.Range("M137").Value = 24.68190000000004
$ws.Range("N137").Value = -10395.2502

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3161.8386
$ws.Range("I138").Value = 2152.5667
$ws.Range("J138").Value = 4108.0312
$ws.Range("K138").Value = 6457.7001
$ws.Range("L138").Value = 12324.0936
$ws.Range("M138").Value = -1317.7001
$ws.Range("N138").Value = -22604.0936

# Row 139: Something Salty and Ceremonial
$ws.Range("H139").Value = 70470
$ws.Range("J139").Value = 70470
$ws.Range("L139").Value = 70470
$ws.Range("N139").Value = -80750

# Row 140: Tome for Tradition
$ws.Range("H140").Value = 73014.28999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6113.7
$ws.Range("I32").Value = 5454.953
$ws.Range("J32").Value = 9846.6
$ws.Range("K32").Value = 5454.953
$ws.Range("L32").Value = 9846.6
$ws.Range("M32").Value = -5167.953
$ws.Range("N32").Value = -10420.6

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1023.5833
$ws.Range("I74").Value = 676.9
$ws.Range("J74").Value = 2757
$ws.Range("K74").Value = 676.9
$ws.Range("L74").Value = 2757
$ws.Range("M74").Value = 197.1
$ws.Range("N74").Value = -4505

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1023.5833
$ws.Range("I77").Value = 676.9
$ws.Range("J77").Value = 2757
$ws.Range("K77").Value = 3384.5
$ws.Range("L77").Value = 13785
$ws.Range("M77").Value = 983.5
$ws.Range("N77").Value = -22521

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 6094.4
$ws.Range("I122").Value = 6960.5386
$ws.Range("J122").Value = 4485.857
$ws.Range("K122").Value = 20881.6158
$ws.Range("L122").Value = 13457.571
$ws.Range("M122").Value = -18431.6158
$ws.Range("N122").Value = -18357.571

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2759.8823
$ws.Range("I132").Value = 2313.4583
$ws.Range("J132").Value = 3831.3
$ws.Range("K132").Value = 6940.374899999999
$ws.Range("L132").Value = 11493.9
$ws.Range("M132").Value = -4410.374899999999
$ws.Range("N132").Value = -16553.9

# Row 138: Don't Ask about the Rivets
$ws.Range("H138").Value = 70016.664
$ws.Range("J138").Value = 70016.664
$ws.Range("L138").Value = 70016.664
$ws.Range("N138").Value = -80296.664

# Row 141: Essays on Equipment
$ws.Range("H141").Value = 61323.08
$ws.Range("I141").Value = 22000
$ws.Range("K141").Value = 22000
$ws.Range("M141").Value = -16820

$ws = $wb.Worksheets.Item("BSM")
# Row 140: Ceremonial Teeth
$ws.Range("H140").Value = 89450
$ws.Range("J140").Value = 89450
$ws.Range("L140").Value = 89450
$ws.Range("N140").Value = -99810

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 8844.465
$ws.Range("I31").Value = 2664.0444
$ws.Range("J31").Value = 15627.854
$ws.Range("K31").Value = 2664.0444
$ws.Range("L31").Value = 15627.854
$ws.Range("M31").Value = -2369.0444
$ws.Range("N31").Value = -16217.854

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 8844.465
$ws.Range("I34").Value = 2664.0444
$ws.Range("J34").Value = 15627.854
$ws.Range("K34").Value = 2664.0444
$ws.Range("L34").Value = 15627.854
$ws.Range("M34").Value = -2462.0444
$ws.Range("N34").Value = -16031.854

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1347.7894
$ws.Range("I58").Value = 740
$ws.Range("J58").Value = 2023.1111
$ws.Range("K58").Value = 740
$ws.Range("L58").Value = 2023.1111
$ws.Range("M58").Value = -537
$ws.Range("N58").Value = -2429.1111

# Row 74: License to Heal
$ws.Range("H74").Value = 16145.556
$ws.Range("J74").Value = 16145.556
$ws.Range("L74").Value = 16145.556
$ws.Range("N74").Value = -17893.556

# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 16145.556
$ws.Range("J77").Value = 16145.556
$ws.Range("L77").Value = 48436.66800000001
$ws.Range("N77").Value = -57172.66800000001

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1147.2
$ws.Range("I132").Value = 830.93335
$ws.Range("J132").Value = 1621.6
$ws.Range("K132").Value = 2492.80005
$ws.Range("L132").Value = 4864.799999999999
$ws.Range("M132").Value = 37.19995000000017
$ws.Range("N132").Value = -9924.799999999999

# Row 136: Turali Quality
$ws.Range("H136").Value = 1347.7894
$ws.Range("I136").Value = 740
$ws.Range("J136").Value = 2023.1111
$ws.Range("K136").Value = 2220
$ws.Range("L136").Value = 6069.3333
$ws.Range("M136").Value = 330
$ws.Range("N136").Value = -11169.3333

# Row 138: Bow Out
$ws.Range("H138").Value = 49914.285
$ws.Range("J138").Value = 49914.285
$ws.Range("L138").Value = 49914.285
$ws.Range("N138").Value = -60194.285

# Row 140: Spear Pressure
$ws.Range("H140").Value = 72333.336
$ws.Range("J140").Value = 72333.336
$ws.Range("L140").Value = 72333.336
$ws.Range("N140").Value = -82693.336

$ws = $wb.Worksheets.Item("CUL")
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 820.17
$ws.Range("I113").Value = 546
$ws.Range("J113").Value = 844.01086
$ws.Range("K113").Value = 1638
$ws.Range("L113").Value = 2532.03258
$ws.Range("M113").Value = 532
$ws.Range("N113").Value = -6872.03258

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 5264006.5
$ws.Range("I131").Value = 918.7619
$ws.Range("J131").Value = 6757585.5
$ws.Range("K131").Value = 2756.2857
$ws.Range("L131").Value = 20272756.5
$ws.Range("M131").Value = 2283.7143
$ws.Range("N131").Value = -20282836.5

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 2045.3914
$ws.Range("I132").Value = 1691.925
$ws.Range("J132").Value = 4401.8335
$ws.Range("K132").Value = 5075.775
$ws.Range("L132").Value = 13205.5005
$ws.Range("M132").Value = -2545.775
$ws.Range("N132").Value = -18265.5005

# Row 138: Orders Anonymous
$ws.Range("H138").Value = 69033.336
$ws.Range("J138").Value = 69033.336
$ws.Range("L138").Value = 69033.336
$ws.Range("N138").Value = -79313.336

# Row 140: The Right Rod
$ws.Range("H140").Value = 89849
$ws.Range("J140").Value = 89849
$ws.Range("L140").Value = 89849
$ws.Range("N140").Value = -100209

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather
$ws.Range("H122").Value = 3281
$ws.Range("I122").Value = 2795.889
$ws.Range("J122").Value = 4008.6667
$ws.Range("K122").Value = 8387.667000000001
$ws.Range("L122").Value = 12026.0001
$ws.Range("M122").Value = -5937.667000000001
$ws.Range("N122").Value = -16926.0001

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 4366.76
$ws.Range("I132").Value = 4395.381
$ws.Range("K132").Value = 13186.143
$ws.Range("M132").Value = -10656.143

# Row 138: Freezing Toes
$ws.Range("H138").Value = 59835.7
$ws.Range("J138").Value = 59835.7
$ws.Range("L138").Value = 59835.7
$ws.Range("N138").Value = -70115.7

# Row 139: Giving Gatherers Their Gear
$ws.Range("H139").Value = 56880
$ws.Range("J139").Value = 69850
$ws.Range("L139").Value = 69850
$ws.Range("N139").Value = -80130

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 5331.4375
$ws.Range("I62").Value = 5725
$ws.Range("J62").Value = 4937.875
$ws.Range("K62").Value = 5725
$ws.Range("L62").Value = 4937.875
$ws.Range("M62").Value = -5101
$ws.Range("N62").Value = -6185.875

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 5331.4375
$ws.Range("I65").Value = 5725
$ws.Range("J65").Value = 4937.875
$ws.Range("K65").Value = 28625
$ws.Range("L65").Value = 24689.375
$ws.Range("M65").Value = -25505
$ws.Range("N65").Value = -30929.375

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1600.4865
$ws.Range("I122").Value = 1190.4333
$ws.Range("K122").Value = 3571.2999
$ws.Range("M122").Value = -1121.2999

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1439.0952
$ws.Range("I126").Value = 1160.7059
$ws.Range("J126").Value = 2622.25
$ws.Range("K126").Value = 3482.1177
$ws.Range("L126").Value = 7866.75
$ws.Range("M126").Value = -1012.1177
$ws.Range("N126").Value = -12806.75

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1030.8334
$ws.Range("I136").Value = 771.8
$ws.Range("J136").Value = 1411.7646
$ws.Range("K136").Value = 2315.4
$ws.Range("L136").Value = 4235.293799999999
$ws.Range("M136").Value = 234.6000000000004
$ws.Range("N136").Value = -9335.293799999999

# Row 138: Halfgloves, Full Effort
$ws.Range("H138").Value = 68366.664
$ws.Range("J138").Value = 68366.664
$ws.Range("L138").Value = 68366.664
$ws.Range("N138").Value = -78646.664

# Row 139: Cruel Climates
$ws.Range("H139").Value = 59840
$ws.Range("J139").Value = 59840
$ws.Range("L139").Value = 59840
$ws.Range("N139").Value = -70120

# Row 140: Glamorous Gloves
$ws.Range("H140").Value = 29950
$ws.Range("J140").Value = 29950
$ws.Range("L140").Value = 29950
$ws.Range("N140").Value = -40310
